$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 255, shifting existing rows 255..323 down to 256..324.
$ws.Rows("255:255").Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(255, 1).Value  = 11
$ws.Cells.Item(255, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(255, 3).Value  = "Bíobío"
$ws.Cells.Item(255, 4).Value  = 44785
$ws.Cells.Item(255, 5).Value  = 8
$ws.Cells.Item(255, 6).Value  = "Fruta"
$ws.Cells.Item(255, 7).Value  = 100102
$ws.Cells.Item(255, 8).Value  = "Cítricos"
$ws.Cells.Item(255, 9).Value  = 100102005
$ws.Cells.Item(255, 10).Value = "Naranja"
$ws.Cells.Item(255, 11).Value = "Lane Late"
$ws.Cells.Item(255, 12).Value = "Primera"
$ws.Cells.Item(255, 13).Value = 220
$ws.Cells.Item(255, 14).Value = 6000
$ws.Cells.Item(255, 15).Value = 6500
$ws.Cells.Item(255, 16).Value = 6273
$ws.Cells.Item(255, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(255, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(255, 19).Value = 418
$ws.Cells.Item(255, 20).Value = 15
